# repull data, push all data, mean calculation
# Updates the dSF (column F) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -6
    9  = -7
    13 = -1
    22 = -3
    26 = -2
    27 = -4
    32 = 3
    33 = 0
    35 = -10
    36 = 3
    37 = 2
    39 = -7
    40 = -5
    42 = -1
    45 = -8
    47 = 2
    50 = 1
    52 = -1
    54 = 5
    58 = 9
    59 = -3
    61 = -4
    62 = 0
    68 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
